$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Sheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 134.27272
$ws.Range("I9").Value = 147.44444
$ws.Range("J9").Value = 75
$ws.Range("K9").Value = 147.44444
$ws.Range("L9").Value = 75
$ws.Range("M9").Value = 21.55556000000001
$ws.Range("N9").Value = -413
# Row 28
$ws.Range("H28").Value = 773.6875
$ws.Range("I28").Value = 773.6875
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 773.6875
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -288.6875
$ws.Range("N28").ClearContents()
# Row 41
$ws.Range("H41").Value = 290.6
$ws.Range("I41").Value = 148.33333
$ws.Range("J41").Value = 504
$ws.Range("K41").Value = 148.33333
$ws.Range("L41").Value = 504
$ws.Range("M41").Value = 291.66667
$ws.Range("N41").Value = -1384
# Row 112
$ws.Range("H112").Value = 1334.6809
$ws.Range("J112").Value = 1334.6809
$ws.Range("L112").Value = 4004.0427
$ws.Range("N112").Value = -6220.0427
# Row 129
$ws.Range("H129").Value = 1336.4857
$ws.Range("I129").Value = 777.25
$ws.Range("J129").Value = 1370.3788
$ws.Range("K129").Value = 2331.75
$ws.Range("L129").Value = 4111.136399999999
$ws.Range("M129").Value = 2668.25
$ws.Range("N129").Value = -14111.1364
# Row 132
$ws.Range("H132").Value = 377992.28
$ws.Range("I132").Value = 164816.67
$ws.Range("J132").Value = 771547.25
$ws.Range("K132").Value = 494450.01
$ws.Range("L132").Value = 2314641.75
$ws.Range("M132").Value = -491920.01
$ws.Range("N132").Value = -2319701.75
# Row 137
$ws.Range("H137").Value = 735868.4
$ws.Range("I137").Value = 2510763.8
$ws.Range("J137").Value = 2759.4565
$ws.Range("K137").Value = 7532291.399999999
$ws.Range("L137").Value = 8278.369499999999
$ws.Range("M137").Value = -7529741.399999999
$ws.Range("N137").Value = -13378.3695
# Row 138
$ws.Range("H138").Value = 2862.3472
$ws.Range("I138").Value = 1356.3572
$ws.Range("J138").Value = 3820.7046
$ws.Range("K138").Value = 4069.0716
$ws.Range("L138").Value = 11462.1138
$ws.Range("M138").Value = 1070.9284
$ws.Range("N138").Value = -21742.1138

# ---- Sheet ARM ----
$ws = $wb.Sheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1888.24
$ws.Range("I32").Value = 1371.2338
$ws.Range("K32").Value = 1371.2338
$ws.Range("M32").Value = -1084.2338
# Row 45
$ws.Range("H45").Value = 2109
$ws.Range("I45").Value = 3012.2
$ws.Range("J45").Value = 1205.8
$ws.Range("K45").Value = 3012.2
$ws.Range("L45").Value = 1205.8
$ws.Range("M45").Value = -2635.2
$ws.Range("N45").Value = -1959.8
# Row 122
$ws.Range("H122").Value = 3248.525
$ws.Range("I122").Value = 2994.7334
$ws.Range("J122").Value = 4009.9
$ws.Range("K122").Value = 8984.200199999999
$ws.Range("L122").Value = 12029.7
$ws.Range("M122").Value = -6534.200199999999
$ws.Range("N122").Value = -16929.7
# Row 132
$ws.Range("H132").Value = 2123.8667
$ws.Range("I132").Value = 1029.4375
$ws.Range("J132").Value = 3374.6428
$ws.Range("K132").Value = 3088.3125
$ws.Range("L132").Value = 10123.9284
$ws.Range("M132").Value = -558.3125
$ws.Range("N132").Value = -15183.9284
# Row 137
$ws.Range("H137").Value = 42115.43
$ws.Range("J137").Value = 42115.43
$ws.Range("L137").Value = 42115.43
$ws.Range("N137").Value = -52315.43

# ---- Sheet BSM ----
$ws = $wb.Sheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 7253.9
$ws.Range("I20").Value = 2688.8235
$ws.Range("J20").Value = 13223.615
$ws.Range("K20").Value = 2688.8235
$ws.Range("L20").Value = 13223.615
$ws.Range("M20").Value = -2441.8235
$ws.Range("N20").Value = -13717.615
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 86
$ws.Range("H86").Value = 2825.75
$ws.Range("I86").Value = 2825.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2825.75
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1702.75
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 2825.75
$ws.Range("I89").Value = 2825.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 14128.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -8512.75
$ws.Range("N89").ClearContents()
# Row 99
$ws.Range("H99").Value = 3030
$ws.Range("I99").Value = 1450
$ws.Range("J99").Value = 5795
$ws.Range("K99").Value = 1450
$ws.Range("L99").Value = 5795
$ws.Range("N99").Value = -8791
$ws.Range("M99").Value = 48

# ---- Sheet CRP ----
$ws = $wb.Sheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 206509.36
$ws.Range("I31").Value = 423119.75
$ws.Range("J31").Value = 2640.7646
$ws.Range("K31").Value = 423119.75
$ws.Range("L31").Value = 2640.7646
$ws.Range("M31").Value = -422824.75
$ws.Range("N31").Value = -3230.7646
# Row 34
$ws.Range("H34").Value = 206509.36
$ws.Range("I34").Value = 423119.75
$ws.Range("J34").Value = 2640.7646
$ws.Range("K34").Value = 423119.75
$ws.Range("L34").Value = 2640.7646
$ws.Range("M34").Value = -422917.75
$ws.Range("N34").Value = -3044.7646
# Row 58
$ws.Range("H58").Value = 2691.6177
$ws.Range("I58").Value = 1503.7407
$ws.Range("K58").Value = 1503.7407
$ws.Range("M58").Value = -1300.7407
# Row 132
$ws.Range("H132").Value = 2474.2927
$ws.Range("I132").Value = 2034.6666
$ws.Range("K132").Value = 6103.9998
$ws.Range("M132").Value = -3573.9998
# Row 134
$ws.Range("H134").Value = 1184.6857
$ws.Range("I134").Value = 795.3125
$ws.Range("J134").Value = 5338
$ws.Range("K134").Value = 2385.9375
$ws.Range("L134").Value = 16014
$ws.Range("M134").Value = 149.0625
$ws.Range("N134").Value = -21084
# Row 136
$ws.Range("H136").Value = 2691.6177
$ws.Range("I136").Value = 1503.7407
$ws.Range("K136").Value = 4511.2221
$ws.Range("M136").Value = -1961.2221

# ---- Sheet CUL ----
$ws = $wb.Sheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1578.1111
$ws.Range("I68").Value = 1299.5385
$ws.Range("K68").Value = 3898.6155
$ws.Range("M68").Value = -3087.6155
# Row 71
$ws.Range("H71").Value = 1578.1111
$ws.Range("I71").Value = 1299.5385
$ws.Range("K71").Value = 11695.8465
$ws.Range("M71").Value = -7639.846500000001
# Row 131
$ws.Range("H131").Value = 805.38
$ws.Range("J131").Value = 828.71576
$ws.Range("L131").Value = 2486.14728
$ws.Range("N131").Value = -12566.14728

# ---- Sheet GSM ----
$ws = $wb.Sheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1179.3636
$ws.Range("I97").Value = 1118.5714
$ws.Range("K97").Value = 1118.5714
$ws.Range("M97").Value = -622.5714
# Row 102
$ws.Range("H102").Value = 2527.257
$ws.Range("I102").Value = 1775.7037
$ws.Range("K102").Value = 1775.7037
$ws.Range("M102").Value = -153.7037
# Row 122
$ws.Range("H122").Value = 4812.5557
$ws.Range("I122").Value = 3330.4285
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 9991.2855
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -7541.2855
$ws.Range("N122").Value = -34900
# Row 126
$ws.Range("H126").Value = 3413.21
$ws.Range("I126").Value = 2782.4084
$ws.Range("J126").Value = 4957.5864
$ws.Range("K126").Value = 8347.225199999999
$ws.Range("L126").Value = 14872.7592
$ws.Range("M126").Value = -5877.225199999999
$ws.Range("N126").Value = -19812.7592
# Row 132
$ws.Range("H132").Value = 2718.761
$ws.Range("I132").Value = 1550.1428
$ws.Range("J132").Value = 4536.6113
$ws.Range("K132").Value = 4650.428400000001
$ws.Range("L132").Value = 13609.8339
$ws.Range("M132").Value = -2120.428400000001
$ws.Range("N132").Value = -18669.8339
# Row 139
$ws.Range("H139").Value = 49999
$ws.Range("J139").Value = 49999
$ws.Range("L139").Value = 49999
$ws.Range("N139").Value = -60279

# ---- Sheet LTW ----
$ws = $wb.Sheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 138.38461
$ws.Range("I55").Value = 130
$ws.Range("J55").Value = 148.16667
$ws.Range("K55").Value = 130
$ws.Range("L55").Value = 148.16667
$ws.Range("M55").Value = 43
$ws.Range("N55").Value = -494.16667
# Row 110
$ws.Range("H110").Value = 29750
$ws.Range("J110").Value = 29750
$ws.Range("L110").Value = 29750
$ws.Range("N110").Value = -37930
# Row 122
$ws.Range("H122").Value = 4397.64
$ws.Range("I122").Value = 3854.3333
$ws.Range("J122").Value = 7250
$ws.Range("K122").Value = 11562.9999
$ws.Range("L122").Value = 21750
$ws.Range("M122").Value = -9112.999899999999
$ws.Range("N122").Value = -26650
# Row 132
$ws.Range("H132").Value = 3997.8484
$ws.Range("I132").Value = 3178.1
$ws.Range("J132").Value = 5259
$ws.Range("K132").Value = 9534.299999999999
$ws.Range("L132").Value = 15777
$ws.Range("M132").Value = -7004.299999999999
$ws.Range("N132").Value = -20837

# ---- Sheet WVR ----
$ws = $wb.Sheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 54438560
$ws.Range("I96").Value = 91864480
$ws.Range("K96").Value = 91864480
$ws.Range("M96").Value = -91863107
# Row 122
$ws.Range("H122").Value = 3855.5518
$ws.Range("I122").Value = 2468.9
$ws.Range("J122").Value = 4585.3687
$ws.Range("K122").Value = 7406.700000000001
$ws.Range("L122").Value = 13756.1061
$ws.Range("M122").Value = -4956.700000000001
$ws.Range("N122").Value = -18656.1061
# Row 123
$ws.Range("H123").Value = 39562.5
$ws.Range("J123").Value = 39562.5
$ws.Range("L123").Value = 39562.5
$ws.Range("N123").Value = -49362.5
# Row 132
$ws.Range("H132").Value = 12348300
$ws.Range("I132").Value = 2198.7273
$ws.Range("J132").Value = 20836246
$ws.Range("K132").Value = 6596.1819
$ws.Range("L132").Value = 62508738
$ws.Range("M132").Value = -4066.1819
$ws.Range("N132").Value = -62513798
